$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be
# auto-converted to numbers by Excel (values are stored as text)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.427.64"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3
$ws.Range("D3").Value = "1.876.51"
$ws.Range("E3").Value = "  +1.01%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "0.7129"
$ws.Range("E5").Value = "  +1.49%  "

# Row 6
$ws.Range("D6").Value = "241.68"
$ws.Range("E6").Value = "  +1.52%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "0.07838"
$ws.Range("E8").Value = "  -2.27%  "

# Row 9
$ws.Range("D9").Value = "0.3112"
$ws.Range("E9").Value = "  +3.00%  "

# Row 10
$ws.Range("D10").Value = "25.15"
$ws.Range("E10").Value = "  +6.95%  "

# Row 11
$ws.Range("D11").Value = "0.08254"
$ws.Range("E11").Value = "  +0.84%  "

# Row 12
$ws.Range("D12").Value = "0.7291"
$ws.Range("E12").Value = "  +3.09%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.261"
$ws.Range("E13").Value = "  +1.24%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.848.51"
$ws.Range("E14").Value = "  -0.76%  "

# Row 15
$ws.Range("D15").Value = "90.96"
$ws.Range("E15").Value = "  +1.55%  "

# Row 16
$ws.Range("D16").Value = "29.438.43"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17
$ws.Range("D17").Value = "5.907"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18
$ws.Range("D18").Value = "246.99"
$ws.Range("E18").Value = "  +3.77%  "

# Row 19
$ws.Range("D19").Value = "0.000007866"
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").Value = "7.961"
$ws.Range("E22").Value = "  +6.62%  "

# Row 23
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("E24").Value = "  +10.25%  "

# Row 25
$ws.Range("D25").Value = "163.75"
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
$ws.Range("D26").Value = "8.992"
$ws.Range("E26").Value = "  +1.15%  "

# Row 27
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  +0.86%  "

# Row 28
$ws.Range("D28").Value = "1.361"
$ws.Range("E28").Value = "  -3.96%  "

# Row 29
$ws.Range("D29").Value = "1.496"
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
$ws.Range("D30").Value = "4.360"
$ws.Range("E30").Value = "  -0.36%  "

# Row 31
$ws.Range("D31").Value = "4.121"
$ws.Range("E31").Value = "  +2.38%  "

# Row 32
$ws.Range("D32").Value = "0.05309"
$ws.Range("E32").Value = "  +2.07%  "

# Row 33
$ws.Range("D33").Value = "1.926"
$ws.Range("E33").Value = "  +0.21%  "

# Row 34
$ws.Range("D34").Value = "1.199"
$ws.Range("E34").Value = "  +3.21%  "

# Row 35
$ws.Range("D35").Value = "0.7223"
$ws.Range("E35").Value = "  +0.61%  "

# Row 36
$ws.Range("D36").Value = "2.681"
$ws.Range("E36").Value = "  -0.22%  "

# Row 37
$ws.Range("D37").Value = "0.01863"
$ws.Range("E37").Value = "  +0.58%  "

# Row 38
$ws.Range("D38").Value = "1.253.92"
$ws.Range("E38").Value = "  +8.80%  "

# Row 39
$ws.Range("D39").Value = "2.730"
$ws.Range("E39").Value = "  +0.14%  "

# Row 40
$ws.Range("D40").Value = "0.9081"
$ws.Range("E40").Value = "  -3.78%  "

# Row 41
$ws.Range("D41").Value = "73.65"
$ws.Range("E41").Value = "  +4.18%  "

# Row 42
$ws.Range("D42").Value = "6.141"
$ws.Range("E42").Value = "  +2.34%  "

# Row 43
$ws.Range("E43").Value = "  +0.08%  "

# Row 44
$ws.Range("D44").Value = "103.40"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45
$ws.Range("E45").Value = "  +0.63%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.012.61"
$ws.Range("E46").Value = "  -1.73%  "

# Row 47
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D47").Value = "2.929"
$ws.Range("E47").Value = "  +12.70%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.766"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("E49").Value = "  +1.80%  "

# Row 50
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.4317"
$ws.Range("E50").Value = "  +1.20%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.237"
$ws.Range("E51").Value = "  +0.82%  "
